$wb = $excel.ActiveWorkbook

# The "展览" (sheet 1) and "全部类型" (sheet 4) sheets both list the same
# animation-convention events. The oldest event (2024-08-10, row 2) has
# passed and is removed from the list; the remaining events shift up one
# row, and the "want to go" counter for the now-first event is bumped
# from 995 to 996 (everything else about that row is unchanged).

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the obsolete first event row; remaining rows shift up,
    # updating the sheet's used range/dimension automatically.
    $ws.Rows.Item(2).Delete()

    # The index numbers in column A (1, 2, 3 for rows 2-4) are not part of
    # the row-shift; they stay fixed even though the rest of each row
    # shifted up, so restore them after the deletion.
    $ws.Range("A2").Value = 1
    $ws.Range("A3").Value = 2
    $ws.Range("A4").Value = 3

    # Bump the "想去人数" (want-to-go count) for the event that is now in
    # row 2 (previously row 3) from 995 to 996.
    $ws.Range("F2").Value = 996
}
